$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A66").Value = "2019年5月22日20:42:42"
$ws.Range("B66").Value = "周三"
$ws.Range("C66").Value = "浅入了解Ajax Json Restful"
$ws.Range("D66").Value = "08:30--10:10"

$ws.Range("C67").Value = "项目包名规范化"
$ws.Range("D67").Value = "18:00--20:30"

$ws.Range("D67").Select()
